$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.946.96'
$ws.Range("E2").Value = '  -0.06%  '

$ws.Range("D3").Value = '2.367.31'
$ws.Range("E3").Value = '  -1.42%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''319.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.75%  '

$ws.Range("D6").Value = '''107.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.01%  '

$ws.Range("D7").Value = '''0.638'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.01%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '''0.625'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.05%  '

$ws.Range("D10").Value = '''41.36'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.85%  '

$ws.Range("D11").Value = '''0.0931'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("D12").Value = '''8.60'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.77%  '

$ws.Range("E13").Value = '  -2.59%  '

$ws.Range("E14").Value = '  +0.27%  '

$ws.Range("D15").Value = '''15.99'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.52%  '

$ws.Range("D16").Value = '2.724.25'
$ws.Range("E16").Value = '  -1.29%  '

$ws.Range("D17").Value = '2.317.52'
$ws.Range("E17").Value = '  -3.49%  '

$ws.Range("D18").Value = '42.908.26'
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").Value = '''7.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.81%  '

$ws.Range("E20").Value = '  -0.36%  '

$ws.Range("D21").Value = '''76.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("E22").Value = '  -6.12%  '

$ws.Range("D23").Value = '''267.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.90%  '

$ws.Range("D24").Value = '''2.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.89%  '

$ws.Range("D25").Value = '''9.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.72%  '

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").Value = '''11.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.55%  '

$ws.Range("D28").Value = '''23.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.23%  '

$ws.Range("E29").Value = '  +2.32%  '

$ws.Range("D30").Value = '''36.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.00%  '

$ws.Range("D31").Value = '''168.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.52%  '

$ws.Range("D32").Value = '''0.0911'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.09%  '

$ws.Range("D33").Value = '''6.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.36%  '

$ws.Range("E34").Value = '  -6.80%  '

$ws.Range("D35").Value = '''0.132'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.81%  '

$ws.Range("E36").Value = '  +11.92%  '

$ws.Range("D37").Value = '''4.77'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("D38").Value = '''0.0364'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.54%  '

$ws.Range("E39").Value = '  -1.68%  '

$ws.Range("D40").Value = '''2.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.20%  '

$ws.Range("D41").Value = '''105.49'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.57%  '

$ws.Range("E42").Value = '  -1.84%  '

$ws.Range("E43").Value = '  +3.25%  '

$ws.Range("D44").Value = '''71.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.52%  '

$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").Value = '''12.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.76%  '

$ws.Range("D47").Value = '''114.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.75%  '

$ws.Range("D48").Value = '''5.56'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.61%  '

$ws.Range("D49").Value = '''9.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.43%  '

$ws.Range("D50").Value = '''76.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.07%  '

$ws.Range("E51").Value = '  +1.49%  '
